# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> "Integral" / "Red Violet"  (bound to the slide master)
#   ppt/theme/theme2.xml -> "Office Theme" / "Office"   (bound to the notes master)
#
# The authored edit swaps the two themes' contents: the slide master's
# theme becomes the "Office" palette (previously on the notes master),
# and the notes master's theme becomes the "Red Violet" palette
# (previously on the slide master). Font scheme / format scheme are
# identical between the two themes, so the only substantive change is
# the 12-slot colour scheme (and the cosmetic theme / colour-scheme
# names) carried by each theme part.
#
# Helper: turn an RRGGBB hex string into the decimal BGR integer that
# PowerPoint's ColorFormat/ThemeColor .RGB property expects (the COM
# RGB() macro encodes R + G*256 + B*65536).
function ConvertTo-RGBValue([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- Target "Office" colour scheme (what theme1.xml should become) ---
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$slideMaster = $p.SlideMaster
$slideMaster.Theme.Name = "Office Theme"

$slideColorScheme = $slideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $slideColorScheme.Item($i).RGB = ConvertTo-RGBValue $officeColors[$i - 1]
}
